# edit.ps1
# Applies the dialogue edits described by the commit ("write some new for stephen")
# to "Act 2 Prim / Scene 2B": new Asher reaction-line stubs are inserted, one of
# his existing lines gets a new emotion tag, the mystery "???" speaker tag for
# Prim's entrance is split into a proper Asher reaction + a Prim line, a split
# run is normalized, the scene-heading text changes, and a missing comma is
# added.

$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new empty "Asher (neutral curious):" line right before his
#    existing "Yeah." line.
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Asher (neutral neutral): Yeah.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$idx = $rng.Paragraphs(1).Index
$rng.InsertParagraphBefore()
$d.Paragraphs($idx).Range.Text = "Asher (neutral curious):"

# ------------------------------------------------------------------
# 2) That "Yeah." line itself now carries a different reaction tag.
# ------------------------------------------------------------------
$d.Content.Find.Execute("Asher (neutral neutral): Yeah.", $true, $false, $false, $false, $false, $true, 1, $false, "Asher (neutral smiling_nervous): Yeah.", 2) | Out-Null

# ------------------------------------------------------------------
# 3) The anonymous "???: Um..." line becomes an Asher reaction stub
#    followed by Prim's actual line.
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("???: Um...", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$idx = $rng.Paragraphs(1).Index
$p = $d.Paragraphs($idx)
$p.Range.Text = "Asher (neutral curious):"
$p.Range.InsertParagraphAfter()
$d.Paragraphs($idx + 1).Range.Text = "?Prim: Um" + [char]0x2026

# ------------------------------------------------------------------
# 4) Normalize the split "She "/"glances"/" nervously at Asher." run into a
#    single run, then add a new "Asher (neutral surprise):" stub after it.
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("She glances nervously at Asher.", $true, $false, $false, $false, $false, $true, 1, $false, "She glances nervously at Asher.", 2) | Out-Null
$rng2 = $d.Content
$rng2.Find.Execute("She glances nervously at Asher.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$idx = $rng2.Paragraphs(1).Index
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$d.Paragraphs($idx + 1).Range.Text = "Asher (neutral surprise):"

# ------------------------------------------------------------------
# 5) Add a new "Asher (neutral playful):" stub after Pro glances at Asher
#    for his nod of approval.
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("I glance at Asher, and he gives me an approving nod.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$idx = $rng.Paragraphs(1).Index
$d.Paragraphs($idx).Range.InsertParagraphAfter()
$d.Paragraphs($idx + 1).Range.Text = "Asher (neutral playful):"

# ------------------------------------------------------------------
# 6) Scene heading rename.
# ------------------------------------------------------------------
$d.Content.Find.Execute("School Grounds", $true, $false, $false, $false, $false, $true, 1, $false, "Back of School", 2) | Out-Null

# ------------------------------------------------------------------
# 7) Missing comma in the narration.
# ------------------------------------------------------------------
$d.Content.Find.Execute("nobody" + [char]0x2019 + "s around Prim lets out", $true, $false, $false, $false, $false, $true, 1, $false, "nobody" + [char]0x2019 + "s around, Prim lets out", 2) | Out-Null
